$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Kitangani" landing site (row 24-25, Takaungu BMU) gains two more
# enumerators: Ngala and Edward Yaa. Insert two new rows right after the
# existing Kitangani rows (before the old row 26 "Vitanga viwili" block),
# which pushes the remaining rows down by two.
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value = "Takaungu"
$ws.Range("B26").Value = "Kitangani"
$ws.Range("C26").Value = "Ngala"

$ws.Range("A27").Value = "Takaungu"
$ws.Range("B27").Value = "Kitangani"
$ws.Range("C27").Value = "Edward Yaa"

# Reflect the updated selection/scroll state from the edit.
[void]$ws.Range("D25").Select()

Write-Host "Inserted enumerators Ngala and Edward Yaa for Kitangani landing site"
